# Injuries_Master_Clubs.xlsx - new scrape pass (2025-11-12)
# 1) Broadhurst Alex (AMR) recovered -> remove his row from "snapshot" and log him on "returned"
# 2) Refresh scraped_at (column K) for every remaining snapshot row (new scrape pass)
# 3) Drozdov Ivan (CSK) is no longer a brand-new entry this run -> clear "new_injured"

$wb = $excel.ActiveWorkbook

$snapshot = $wb.Worksheets.Item("snapshot")
$returned = $wb.Worksheets.Item("returned")
$newInjured = $wb.Worksheets.Item("new_injured")

# --- 1) Remove the Broadhurst Alex row (row 10) from the snapshot ---
$snapshot.Rows.Item(10).Delete()

# --- 2) Refresh the scraped_at column (K) for the remaining 51 rows (2..52) ---
$newScrapedAt = @(
    '2025-11-12T09:23:16.117744+00:00',
    '2025-11-12T09:23:16.117768+00:00',
    '2025-11-12T09:23:16.117781+00:00',
    '2025-11-12T09:23:18.816979+00:00',
    '2025-11-12T09:23:18.817023+00:00',
    '2025-11-12T09:23:18.817036+00:00',
    '2025-11-12T09:23:20.832321+00:00',
    '2025-11-12T09:23:23.281387+00:00',
    '2025-11-12T09:23:23.281421+00:00',
    '2025-11-12T09:23:25.824617+00:00',
    '2025-11-12T09:23:25.824647+00:00',
    '2025-11-12T09:23:25.824665+00:00',
    '2025-11-12T09:23:25.824682+00:00',
    '2025-11-12T09:23:30.114612+00:00',
    '2025-11-12T09:23:32.012341+00:00',
    '2025-11-12T09:23:33.975310+00:00',
    '2025-11-12T09:23:33.975337+00:00',
    '2025-11-12T09:23:33.975359+00:00',
    '2025-11-12T09:23:36.309807+00:00',
    '2025-11-12T09:23:39.395462+00:00',
    '2025-11-12T09:23:39.395491+00:00',
    '2025-11-12T09:23:41.308603+00:00',
    '2025-11-12T09:23:41.308634+00:00',
    '2025-11-12T09:23:41.308653+00:00',
    '2025-11-12T09:23:43.623865+00:00',
    '2025-11-12T09:23:43.623899+00:00',
    '2025-11-12T09:23:43.623920+00:00',
    '2025-11-12T09:23:43.623940+00:00',
    '2025-11-12T09:23:43.623957+00:00',
    '2025-11-12T09:23:45.568999+00:00',
    '2025-11-12T09:23:45.569029+00:00',
    '2025-11-12T09:23:48.492691+00:00',
    '2025-11-12T09:23:48.492726+00:00',
    '2025-11-12T09:23:48.492751+00:00',
    '2025-11-12T09:23:50.808327+00:00',
    '2025-11-12T09:23:50.808363+00:00',
    '2025-11-12T09:23:50.808385+00:00',
    '2025-11-12T09:23:52.703282+00:00',
    '2025-11-12T09:23:52.703312+00:00',
    '2025-11-12T09:23:52.703331+00:00',
    '2025-11-12T09:23:52.703351+00:00',
    '2025-11-12T09:23:52.703369+00:00',
    '2025-11-12T09:23:52.703388+00:00',
    '2025-11-12T09:23:54.549585+00:00',
    '2025-11-12T09:23:54.549617+00:00',
    '2025-11-12T09:23:58.829989+00:00',
    '2025-11-12T09:23:58.830016+00:00',
    '2025-11-12T09:23:58.830034+00:00',
    '2025-11-12T09:23:58.830051+00:00',
    '2025-11-12T09:24:01.169066+00:00',
    '2025-11-12T09:24:01.169097+00:00'
)

for ($i = 0; $i -lt $newScrapedAt.Length; $i++) {
    $row = $i + 2
    $snapshot.Cells.Item($row, 11).Value = $newScrapedAt[$i]
}

# --- 3) Record the recovered player on the "returned" sheet ---
$returned.Cells.Item(2, 1).Value = "АМР"
$returned.Cells.Item(2, 2).Value = "Амур"
$returned.Cells.Item(2, 3).Value = "Броадхёрст Алекс"
$returned.Cells.Item(2, 4).Value = "1369_АМР_броадхерсталекс"
$returned.Cells.Item(2, 5).Value = "RETURN"
$returned.Cells.Item(2, 6).Value = "2025-11-12T17:24:01.672404+08:00"
$returned.Cells.Item(2, 7).NumberFormat = "@"
$returned.Cells.Item(2, 7).Value = "2025-11-12"

# --- 4) Drozdov Ivan is no longer newly-injured -> clear the new_injured data row ---
$newInjured.Range("A2:G2").ClearContents()
